# Applies the "New data structure and gen expression" commit:
#  - expands Sheet1 from 3 rows (2 datasets) to 12 rows (11 datasets), all
#    compared against hPSC
#  - rewrites the header row and every data row's Dataset/Description/
#    Condition1/Condition2/Location columns
#  - widens column A, updates the current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bank the small-print (Arial 10) formatting that currently lives on A3
# (it needs to survive onto the new "Cerebral Organoids d40" row's
# Description cell) in a scratch cell far outside the table before we
# wipe the old 2-dataset rows.
$ws.Range("A3").Copy()
$ws.Range("ZZ100").PasteSpecial(-4122)

# Drop the old table (contents + formatting) and start fresh. Keep the
# range tight to the cells that actually held data/format before so we
# don't spray empty placeholder cells into rows/cols that never had any
# (e.g. F2:G3, which were never part of the old table).
$ws.Range("A1:E3").ClearContents()
$ws.Range("A1:E3").ClearFormats()
$ws.Range("F1:G1").ClearFormats()

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Dataset"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Condition1"
$ws.Range("D1").Value = "Condition2"
$ws.Range("E1").Value = "Location"
$ws.Range("A1:G1").Font.Bold = $true

# ---- Data rows --------------------------------------------------------
$data = @(
  @("NPC P1",                         "NPC at passage 1 in NPM",                                     "NPC P1",                         "hPSC", "data/datasets/NPC_P1_vs_hPSC.Rds"),
  @("NPC P3",                         "NPC at passage 3 in NPM",                                     "NPC P3",                         "hPSC", "data/datasets/NPC_P3_vs_hPSC.Rds"),
  @("NPC P5",                         "NPC at passage 5 in NPM",                                     "NPC P5",                         "hPSC", "data/datasets/NPC_P5_vs_hPSC.Rds"),
  @("Neural Crest d6",                "Neural Crest after day 6 differentiation",                    "Neural Crest d6",                "hPSC", "data/datasets/Neural Crest_d6_vs_hPSC.Rds"),
  @("Sensory Neuron Diff d12",        "Sensory Neurons after 6 additional days of differentiation",  "Sensory Neuron Diff d12",        "hPSC", "data/datasets/Sensory Neuron Diff_d12_vs_hPSC.Rds"),
  @("Sensory Neuron Maturation d18",  "Sensory Neurons after 6 additional days of maturation",       "Sensory Neuron Maturation d18",  "hPSC", "data/datasets/Sensory Neuron Maturation_d18_vs_hPSC.Rds"),
  @("Cerebral Organoids d40",         "Cerebral Organoids day 40 compared to hPSC",                  "Cerebral Organoids d40",         "hPSC", "data/datasets/Cerebral Organoids_d40_vs_hPSC.Rds"),
  @("Dorsal Forebrain d25",           "Dorsal Forebrain Organoids at day 25",                        "Dorsal Forebrain d25",           "hPSC", "data/datasets/Dorsal Forebrain_d25_vs_hPSC.Rds"),
  @("Dorsal Forebrain d50",           "Dorsal Forebrain Organoids at day 50",                        "Dorsal Forebrain d50",           "hPSC", "data/datasets/Dorsal Forebrain_d50_vs_hPSC.Rds"),
  @("Dorsal Forebrain d75",           "Dorsal Forebrain Organoids at day 75",                        "Dorsal Forebrain d75",           "hPSC", "data/datasets/Dorsal Forebrain_d75_vs_hPSC.Rds"),
  @("Ventral Forebrain d25",          "Ventral Forebrain Organoids at day 25",                       "Ventral Forebrain d25",          "hPSC", "data/datasets/Ventral Forebrain_d25_vs_hPSC.Rds")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# The "Cerebral Organoids d40" row (row 8) keeps the old small-print
# (Arial 10) description styling that used to sit on the single
# pre-existing dataset row -- restore it from the scratch cell.
$ws.Range("ZZ100").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("ZZ100").Clear()

# ---- Column width / view state ---------------------------------------
# Target stored width is 32.7109375 characters; the host only resolves
# ColumnWidth onto a 1/6-character pixel grid, so 31.8 is the input that
# lands on the closest representable stored width (32.6666...).
$ws.Columns.Item(1).ColumnWidth = 31.8

$ws.Range("J13").Select()
